# The presentation's main theme (ppt/theme/theme1.xml, bound to the
# slide master) is swapped from the "Integral" / "Red Violet" palette
# to the standard Office Theme palette (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink). The font scheme and format scheme are already
# identical between the two themes, so only the twelve theme colours
# need to move.
#
# Colours are written through Slide.ThemeColorScheme, whose twelve
# slots map - in order - to dk1, lt1, dk2, lt2, accent1, accent2,
# accent3, accent4, accent5, accent6, hlink, folHlink. PowerPoint's
# .RGB values are plain 0xBBGGRR integers (standard OLE COLORREF), so
# each target hex colour below is byte-reversed before being assigned.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette, in ThemeColorScheme slot order 1..12.
$officeThemeRGB = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeThemeRGB[$i - 1]
}
